# Auto-generated Excel COM-interop script
# Applies numeric corrections to several cells across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# as produced by the scheduled profit-recalculation runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 343.5
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H109").Value = 70000
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72774
$ws.Range("H116").Value = 55100.332
$ws.Range("I116").Value = 91327.42999999999
$ws.Range("J116").Value = 4382.4
$ws.Range("K116").Value = 91327.42999999999
$ws.Range("L116").Value = 4382.4
$ws.Range("M116").Value = -87885.42999999999
$ws.Range("N116").Value = -11266.4
$ws.Range("H132").Value = 440458.88
$ws.Range("J132").Value = 6451.2
$ws.Range("L132").Value = 19353.6
$ws.Range("N132").Value = -24413.6
$ws.Range("H135").Value = 1586.96
$ws.Range("I135").Value = 1586.96
$ws.Range("K135").Value = 14282.64
$ws.Range("M135").Value = -11747.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 313.5
$ws.Range("J5").Value = 401.7143
$ws.Range("L5").Value = 401.7143
$ws.Range("N5").Value = -625.7143
$ws.Range("H32").Value = 6045.11
$ws.Range("I32").Value = 6045.11
$ws.Range("K32").Value = 6045.11
$ws.Range("M32").Value = -5758.11
$ws.Range("H61").Value = 7481.7188
$ws.Range("I61").Value = 4776.8
$ws.Range("K61").Value = 4776.8
$ws.Range("M61").Value = -4564.8
$ws.Range("H74").Value = 2720278.2
$ws.Range("I74").Value = 3789561
$ws.Range("K74").Value = 3789561
$ws.Range("M74").Value = -3788687
$ws.Range("H77").Value = 2720278.2
$ws.Range("I77").Value = 3789561
$ws.Range("K77").Value = 18947805
$ws.Range("M77").Value = -18943437
$ws.Range("H94").Value = 12000
$ws.Range("J94").Value = 12000
$ws.Range("L94").Value = 12000
$ws.Range("N94").Value = -13802
$ws.Range("H136").Value = 7481.7188
$ws.Range("I136").Value = 4776.8
$ws.Range("K136").Value = 14330.4
$ws.Range("M136").Value = -11780.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 313.5
$ws.Range("J4").Value = 401.7143
$ws.Range("L4").Value = 401.7143
$ws.Range("N4").Value = -631.7143
$ws.Range("H82").Value = 30146.75
$ws.Range("J82").Value = 65805.71000000001
$ws.Range("L82").Value = 65805.71000000001
$ws.Range("N82").Value = -66571.71000000001
$ws.Range("H85").Value = 30146.75
$ws.Range("J85").Value = 65805.71000000001
$ws.Range("L85").Value = 65805.71000000001
$ws.Range("N85").Value = -68457.71000000001
$ws.Range("H104").Value = 80684
$ws.Range("J104").Value = 80684
$ws.Range("L104").Value = 80684
$ws.Range("N104").Value = -87672
$ws.Range("H107").Value = 889
$ws.Range("I107").Value = 889
$ws.Range("K107").Value = 889
$ws.Range("M107").Value = 1031
$ws.Range("H108").Value = 100000
$ws.Range("J108").Value = 100000
$ws.Range("L108").Value = 100000
$ws.Range("N108").Value = -107680
$ws.Range("H110").Value = 70666.664
$ws.Range("J110").Value = 70666.664
$ws.Range("L110").Value = 70666.664
$ws.Range("N110").Value = -78846.664
$ws.Range("H117").Value = 50742
$ws.Range("J117").Value = 50742
$ws.Range("L117").Value = 50742
$ws.Range("N117").Value = -59920
$ws.Range("H134").Value = 543949.7
$ws.Range("I134").Value = 702179.2
$ws.Range("J134").Value = 5969.5
$ws.Range("K134").Value = 2106537.6
$ws.Range("L134").Value = 17908.5
$ws.Range("N134").Value = -22978.5
$ws.Range("M134").Value = -2104002.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H22").Value = 993
$ws.Range("I22").Value = 993
$ws.Range("K22").Value = 993
$ws.Range("M22").Value = -643
$ws.Range("H31").Value = 8362.885
$ws.Range("I31").Value = 4034.8462
$ws.Range("J31").Value = 12690.923
$ws.Range("K31").Value = 4034.8462
$ws.Range("L31").Value = 12690.923
$ws.Range("M31").Value = -3739.8462
$ws.Range("N31").Value = -13280.923
$ws.Range("H34").Value = 8362.885
$ws.Range("I34").Value = 4034.8462
$ws.Range("J34").Value = 12690.923
$ws.Range("K34").Value = 4034.8462
$ws.Range("L34").Value = 12690.923
$ws.Range("M34").Value = -3832.8462
$ws.Range("N34").Value = -13094.923
$ws.Range("H86").Value = 2812.7856
$ws.Range("I86").Value = 2811.8572
$ws.Range("J86").Value = 2813.7144
$ws.Range("K86").Value = 2811.8572
$ws.Range("L86").Value = 2813.7144
$ws.Range("M86").Value = -1688.8572
$ws.Range("N86").Value = -5059.7144
$ws.Range("H89").Value = 2812.7856
$ws.Range("I89").Value = 2811.8572
$ws.Range("J89").Value = 2813.7144
$ws.Range("K89").Value = 14059.286
$ws.Range("L89").Value = 14068.572
$ws.Range("M89").Value = -8443.286
$ws.Range("N89").Value = -25300.572
$ws.Range("H92").Value = 60000
$ws.Range("J92").Value = 60000
$ws.Range("L92").Value = 60000
$ws.Range("N92").Value = -64992
$ws.Range("H94").Value = 1210.1923
$ws.Range("I94").Value = 869.7692
$ws.Range("J94").Value = 1550.6154
$ws.Range("K94").Value = 869.7692
$ws.Range("L94").Value = 1550.6154
$ws.Range("M94").Value = -418.7692
$ws.Range("N94").Value = -2452.6154
$ws.Range("H108").Value = 69000
$ws.Range("J108").Value = 69000
$ws.Range("L108").Value = 69000
$ws.Range("N108").Value = -76680
$ws.Range("H134").Value = 99440.92
$ws.Range("I134").Value = 105234
$ws.Range("J134").Value = 80130.664
$ws.Range("K134").Value = 315702
$ws.Range("L134").Value = 240391.992
$ws.Range("M134").Value = -313167
$ws.Range("N134").Value = -245461.992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 3161755.2
$ws.Range("I7").Value = 2857236.8
$ws.Range("K7").Value = 8571710.399999999
$ws.Range("M7").Value = -8571598.399999999
$ws.Range("H41").Value = 407.66666
$ws.Range("I41").Value = 82.666664
$ws.Range("J41").Value = 732.6667
$ws.Range("K41").Value = 247.999992
$ws.Range("L41").Value = 2198.0001
$ws.Range("M41").Value = 90.00000800000001
$ws.Range("N41").Value = -2874.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1087.1428
$ws.Range("I2").Value = 74.13333
$ws.Range("K2").Value = 74.13333
$ws.Range("M2").Value = 38.86667
$ws.Range("H70").Value = 6205.375
$ws.Range("I70").Value = 5863.4
$ws.Range("J70").Value = 6775.3335
$ws.Range("K70").Value = 5863.4
$ws.Range("L70").Value = 6775.3335
$ws.Range("M70").Value = -5593.4
$ws.Range("N70").Value = -7315.3335
$ws.Range("H73").Value = 6205.375
$ws.Range("I73").Value = 5863.4
$ws.Range("J73").Value = 6775.3335
$ws.Range("K73").Value = 5863.4
$ws.Range("L73").Value = 6775.3335
$ws.Range("M73").Value = -4927.4
$ws.Range("N73").Value = -8647.333500000001
$ws.Range("H93").Value = 59375
$ws.Range("J93").Value = 59375
$ws.Range("L93").Value = 59375
$ws.Range("N93").Value = -63119
$ws.Range("H97").Value = 2801.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2801.5
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2801.5
$ws.Range("N97").Value = -3793.5
$ws.Range("M97").ClearContents()
$ws.Range("H126").Value = 1114223.6
$ws.Range("I126").Value = 1669366.2
$ws.Range("J126").Value = 3938.6
$ws.Range("K126").Value = 5008098.6
$ws.Range("L126").Value = 11815.8
$ws.Range("M126").Value = -5005628.6
$ws.Range("N126").Value = -16755.8
$ws.Range("H132").Value = 2283.9546
$ws.Range("I132").Value = 1980.919
$ws.Range("K132").Value = 5942.757000000001
$ws.Range("M132").Value = -3412.757000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3908.8572
$ws.Range("I136").Value = 2815.9092
$ws.Range("K136").Value = 8447.7276
$ws.Range("M136").Value = -5897.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2026
$ws.Range("I96").Value = 1589.7778
$ws.Range("J96").Value = 3334.6667
$ws.Range("K96").Value = 1589.7778
$ws.Range("L96").Value = 3334.6667
$ws.Range("M96").Value = -216.7778000000001
$ws.Range("N96").Value = -6080.6667
$ws.Range("H107").Value = 1200.8276
$ws.Range("I107").Value = 963
$ws.Range("K107").Value = 2889
$ws.Range("M107").Value = -969
$ws.Range("H132").Value = 33479998
$ws.Range("I132").Value = 2471603.2
$ws.Range("K132").Value = 7414809.600000001
$ws.Range("M132").Value = -7412279.600000001

Write-Output "Applied 218 cell updates across 8 sheets."
